# تعديل حدث في Card24 - الصف 19
# Moves the service-log entry one row down: row 21's Date/Event/Correction
# (and the row counter in column A) are cleared, and row 22's corresponding
# cells are set to "nan".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Row 21: clear the row counter and the date/event/correction entry
$ws.Range("A21").Value = ""
$ws.Range("L21").Value = ""
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = ""

# Row 22: fill the row counter and the date/event/correction cells with "nan"
$ws.Range("A22").Value = "nan"
$ws.Range("L22").Value = "nan"
$ws.Range("M22").Value = "nan"
$ws.Range("N22").Value = "nan"
